$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1226.9  # H6: 865 -> 1226.9
$ws.Cells.Item(6, 9).Value = 138.42857  # I6: 131.25 -> 138.42857
$ws.Cells.Item(6, 10).Value = 3766.6667  # J6: 3800 -> 3766.6667
$ws.Cells.Item(6, 11).Value = 415.28571  # K6: 393.75 -> 415.28571
$ws.Cells.Item(6, 12).Value = 11300.0001  # L6: 11400 -> 11300.0001
$ws.Cells.Item(6, 13).Value = -303.28571  # M6: -281.75 -> -303.28571
$ws.Cells.Item(6, 14).Value = -11524.0001  # N6: -11624 -> -11524.0001

$ws.Cells.Item(10, 8).Value = 29999.5  # H10: 22668 -> 29999.5
$ws.Cells.Item(10, 10).Value = 29999.5  # J10: 22668 -> 29999.5
$ws.Cells.Item(10, 12).Value = 29999.5  # L10: 22668 -> 29999.5
$ws.Cells.Item(10, 14).Value = -30585.5  # N10: -23254 -> -30585.5

$ws.Cells.Item(17, 8).Value = 680.4578  # H17: 652.0421 -> 680.4578
$ws.Cells.Item(17, 10).Value = 573.2059  # J17: 555.55 -> 573.2059
$ws.Cells.Item(17, 12).Value = 1719.6177  # L17: 1666.65 -> 1719.6177
$ws.Cells.Item(17, 14).Value = -2055.6177  # N17: -2002.65 -> -2055.6177

$ws.Cells.Item(112, 8).Value = 1576.5714  # H112: 1583.4182 -> 1576.5714
$ws.Cells.Item(112, 10).Value = 1593.3334  # J112: 1600.7548 -> 1593.3334
$ws.Cells.Item(112, 12).Value = 4780.0002  # L112: 4802.2644 -> 4780.0002
$ws.Cells.Item(112, 14).Value = -6996.0002  # N112: -7018.2644 -> -6996.0002

$ws.Cells.Item(134, 8).Value = 47865.9  # H134: 48511.105 -> 47865.9
$ws.Cells.Item(134, 9).Value = 29800  # I134: 0 -> 29800
$ws.Cells.Item(134, 11).Value = 29800  # K134: 0 -> 29800
$ws.Cells.Item(134, 13).Value = -24730  # M134: None -> -24730

$ws.Cells.Item(137, 8).Value = 2428.2083  # H137: 2618.432 -> 2428.2083
$ws.Cells.Item(137, 9).Value = 1156.88  # I137: 1269.5 -> 1156.88
$ws.Cells.Item(137, 10).Value = 3810.087  # J137: 3967.3635 -> 3810.087
$ws.Cells.Item(137, 11).Value = 3470.64  # K137: 3808.5 -> 3470.64
$ws.Cells.Item(137, 12).Value = 11430.261  # L137: 11902.0905 -> 11430.261
$ws.Cells.Item(137, 13).Value = -920.6400000000003  # M137: -1258.5 -> -920.6400000000003
$ws.Cells.Item(137, 14).Value = -16530.261  # N137: -17002.0905 -> -16530.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 0  # H17: 25500 -> 0
$ws.Cells.Item(17, 10).Value = 0  # J17: 25500 -> 0
$ws.Cells.Item(17, 12).Value = 0  # L17: 25500 -> 0
$ws.Cells.Item(17, 14).Value = ""  # N17: -25846 -> (removed)

$ws.Cells.Item(61, 8).Value = 1295  # H61: 1419.5 -> 1295
$ws.Cells.Item(61, 9).Value = 980.2308  # I61: 1132.909 -> 980.2308
$ws.Cells.Item(61, 10).Value = 1879.5714  # J61: 1869.8572 -> 1879.5714
$ws.Cells.Item(61, 11).Value = 980.2308  # K61: 1132.909 -> 980.2308
$ws.Cells.Item(61, 12).Value = 1879.5714  # L61: 1869.8572 -> 1879.5714
$ws.Cells.Item(61, 13).Value = -768.2308  # M61: -920.9090000000001 -> -768.2308
$ws.Cells.Item(61, 14).Value = -2303.5714  # N61: -2293.8572 -> -2303.5714

$ws.Cells.Item(136, 8).Value = 1295  # H136: 1419.5 -> 1295
$ws.Cells.Item(136, 9).Value = 980.2308  # I136: 1132.909 -> 980.2308
$ws.Cells.Item(136, 10).Value = 1879.5714  # J136: 1869.8572 -> 1879.5714
$ws.Cells.Item(136, 11).Value = 2940.6924  # K136: 3398.727 -> 2940.6924
$ws.Cells.Item(136, 12).Value = 5638.7142  # L136: 5609.571599999999 -> 5638.7142
$ws.Cells.Item(136, 13).Value = -390.6923999999999  # M136: -848.7270000000003 -> -390.6923999999999
$ws.Cells.Item(136, 14).Value = -10738.7142  # N136: -10709.5716 -> -10738.7142

$ws.Cells.Item(139, 8).Value = 43463.125  # H139: 43369.523 -> 43463.125
$ws.Cells.Item(139, 10).Value = 43463.125  # J139: 43369.523 -> 43463.125
$ws.Cells.Item(139, 12).Value = 43463.125  # L139: 43369.523 -> 43463.125
$ws.Cells.Item(139, 14).Value = -53743.125  # N139: -53649.523 -> -53743.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(69, 8).Value = 0  # H69: 24500 -> 0
$ws.Cells.Item(69, 9).Value = 0  # I69: 14000 -> 0
$ws.Cells.Item(69, 10).Value = 0  # J69: 35000 -> 0
$ws.Cells.Item(69, 11).Value = 0  # K69: 14000 -> 0
$ws.Cells.Item(69, 12).Value = 0  # L69: 35000 -> 0
$ws.Cells.Item(69, 13).Value = ""  # M69: -13189 -> (removed)
$ws.Cells.Item(69, 14).Value = ""  # N69: -36622 -> (removed)

$ws.Cells.Item(72, 8).Value = 0  # H72: 24500 -> 0
$ws.Cells.Item(72, 9).Value = 0  # I72: 14000 -> 0
$ws.Cells.Item(72, 10).Value = 0  # J72: 35000 -> 0
$ws.Cells.Item(72, 11).Value = 0  # K72: 42000 -> 0
$ws.Cells.Item(72, 12).Value = 0  # L72: 105000 -> 0
$ws.Cells.Item(72, 13).Value = ""  # M72: -37944 -> (removed)
$ws.Cells.Item(72, 14).Value = ""  # N72: -113112 -> (removed)

$ws.Cells.Item(86, 8).Value = 1606.25  # H86: 1566.6666 -> 1606.25
$ws.Cells.Item(86, 9).Value = 1550  # I86: 1566.6666 -> 1550
$ws.Cells.Item(86, 10).Value = 2000  # J86: 0 -> 2000
$ws.Cells.Item(86, 11).Value = 1550  # K86: 1566.6666 -> 1550
$ws.Cells.Item(86, 12).Value = 2000  # L86: 0 -> 2000
$ws.Cells.Item(86, 13).Value = -427  # M86: -443.6666 -> -427
$ws.Cells.Item(86, 14).Value = -4246  # N86: None -> -4246

$ws.Cells.Item(89, 8).Value = 1606.25  # H89: 1566.6666 -> 1606.25
$ws.Cells.Item(89, 9).Value = 1550  # I89: 1566.6666 -> 1550
$ws.Cells.Item(89, 10).Value = 2000  # J89: 0 -> 2000
$ws.Cells.Item(89, 11).Value = 7750  # K89: 7833.333000000001 -> 7750
$ws.Cells.Item(89, 12).Value = 10000  # L89: 0 -> 10000
$ws.Cells.Item(89, 13).Value = -2134  # M89: -2217.333000000001 -> -2134
$ws.Cells.Item(89, 14).Value = -21232  # N89: None -> -21232

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(13, 8).Value = 11249.75  # H13: 8999.666999999999 -> 11249.75
$ws.Cells.Item(13, 9).Value = 0  # I13: 1000 -> 0
$ws.Cells.Item(13, 10).Value = 11249.75  # J13: 12999.5 -> 11249.75
$ws.Cells.Item(13, 11).Value = 0  # K13: 1000 -> 0
$ws.Cells.Item(13, 12).Value = 11249.75  # L13: 12999.5 -> 11249.75
$ws.Cells.Item(13, 13).Value = ""  # M13: -861 -> (removed)
$ws.Cells.Item(13, 14).Value = -11527.75  # N13: -13277.5 -> -11527.75

$ws.Cells.Item(22, 8).Value = 509.0909  # H22: 503.60605 -> 509.0909
$ws.Cells.Item(22, 9).Value = 317.95834  # I22: 311.25 -> 317.95834
$ws.Cells.Item(22, 10).Value = 1018.7778  # J22: 1016.55554 -> 1018.7778
$ws.Cells.Item(22, 11).Value = 317.95834  # K22: 311.25 -> 317.95834
$ws.Cells.Item(22, 12).Value = 1018.7778  # L22: 1016.55554 -> 1018.7778
$ws.Cells.Item(22, 13).Value = 32.04165999999998  # M22: 38.75 -> 32.04165999999998
$ws.Cells.Item(22, 14).Value = -1718.7778  # N22: -1716.55554 -> -1718.7778

$ws.Cells.Item(25, 8).Value = 7592.6  # H25: 8985.333000000001 -> 7592.6
$ws.Cells.Item(25, 9).Value = 1000  # I25: 0 -> 1000
$ws.Cells.Item(25, 10).Value = 9240.75  # J25: 8985.333000000001 -> 9240.75
$ws.Cells.Item(25, 11).Value = 1000  # K25: 0 -> 1000
$ws.Cells.Item(25, 12).Value = 9240.75  # L25: 8985.333000000001 -> 9240.75
$ws.Cells.Item(25, 13).Value = -826  # M25: None -> -826
$ws.Cells.Item(25, 14).Value = -9588.75  # N25: -9333.333000000001 -> -9588.75

$ws.Cells.Item(31, 8).Value = 7465165.5  # H31: 7355393.5 -> 7465165.5
$ws.Cells.Item(31, 10).Value = 20837826  # J31: 20004340 -> 20837826
$ws.Cells.Item(31, 12).Value = 20837826  # L31: 20004340 -> 20837826
$ws.Cells.Item(31, 14).Value = -20838416  # N31: -20004930 -> -20838416

$ws.Cells.Item(34, 8).Value = 7465165.5  # H34: 7355393.5 -> 7465165.5
$ws.Cells.Item(34, 10).Value = 20837826  # J34: 20004340 -> 20837826
$ws.Cells.Item(34, 12).Value = 20837826  # L34: 20004340 -> 20837826
$ws.Cells.Item(34, 14).Value = -20838230  # N34: -20004744 -> -20838230

$ws.Cells.Item(58, 8).Value = 1656.4105  # H58: 1742.6067 -> 1656.4105
$ws.Cells.Item(58, 9).Value = 1468.1757  # I58: 1564.3823 -> 1468.1757
$ws.Cells.Item(58, 11).Value = 1468.1757  # K58: 1564.3823 -> 1468.1757
$ws.Cells.Item(58, 13).Value = -1265.1757  # M58: -1361.3823 -> -1265.1757

$ws.Cells.Item(99, 8).Value = 9529197  # H99: 9529207 -> 9529197
$ws.Cells.Item(99, 9).Value = 16670612  # I99: 15388411 -> 16670612
$ws.Cells.Item(99, 10).Value = 7311.5557  # J99: 8000.5 -> 7311.5557
$ws.Cells.Item(99, 11).Value = 16670612  # K99: 15388411 -> 16670612
$ws.Cells.Item(99, 12).Value = 7311.5557  # L99: 8000.5 -> 7311.5557
$ws.Cells.Item(99, 13).Value = -16669114  # M99: -15386913 -> -16669114
$ws.Cells.Item(99, 14).Value = -10307.5557  # N99: -10996.5 -> -10307.5557

$ws.Cells.Item(126, 8).Value = 9529197  # H126: 9529207 -> 9529197
$ws.Cells.Item(126, 9).Value = 16670612  # I126: 15388411 -> 16670612
$ws.Cells.Item(126, 10).Value = 7311.5557  # J126: 8000.5 -> 7311.5557
$ws.Cells.Item(126, 11).Value = 50011836  # K126: 46165233 -> 50011836
$ws.Cells.Item(126, 12).Value = 21934.6671  # L126: 24001.5 -> 21934.6671
$ws.Cells.Item(126, 13).Value = -50009366  # M126: -46162763 -> -50009366
$ws.Cells.Item(126, 14).Value = -26874.6671  # N126: -28941.5 -> -26874.6671

$ws.Cells.Item(134, 8).Value = 3743.9167  # H134: 3944.311 -> 3743.9167
$ws.Cells.Item(134, 9).Value = 5177.7085  # I134: 5371.087 -> 5177.7085
$ws.Cells.Item(134, 10).Value = 2310.125  # J134: 2452.682 -> 2310.125
$ws.Cells.Item(134, 11).Value = 15533.1255  # K134: 16113.261 -> 15533.1255
$ws.Cells.Item(134, 12).Value = 6930.375  # L134: 7358.045999999999 -> 6930.375
$ws.Cells.Item(134, 13).Value = -12998.1255  # M134: -13578.261 -> -12998.1255
$ws.Cells.Item(134, 14).Value = -12000.375  # N134: -12428.046 -> -12000.375

$ws.Cells.Item(136, 8).Value = 1656.4105  # H136: 1742.6067 -> 1656.4105
$ws.Cells.Item(136, 9).Value = 1468.1757  # I136: 1564.3823 -> 1468.1757
$ws.Cells.Item(136, 11).Value = 4404.5271  # K136: 4693.1469 -> 4404.5271
$ws.Cells.Item(136, 13).Value = -1854.5271  # M136: -2143.1469 -> -1854.5271

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1375.0526  # H5: 1435.4722 -> 1375.0526
$ws.Cells.Item(5, 9).Value = 328.6316  # I5: 342.72223 -> 328.6316
$ws.Cells.Item(5, 10).Value = 2421.4736  # J5: 2528.2222 -> 2421.4736
$ws.Cells.Item(5, 11).Value = 985.8948  # K5: 1028.16669 -> 985.8948
$ws.Cells.Item(5, 12).Value = 7264.4208  # L5: 7584.6666 -> 7264.4208
$ws.Cells.Item(5, 13).Value = -873.8948  # M5: -916.16669 -> -873.8948
$ws.Cells.Item(5, 14).Value = -7488.4208  # N5: -7808.6666 -> -7488.4208

$ws.Cells.Item(68, 8).Value = 13446.875  # H68: 13472.125 -> 13446.875
$ws.Cells.Item(68, 9).Value = 690  # I68: 735 -> 690
$ws.Cells.Item(68, 10).Value = 21101  # J68: 17717.834 -> 21101
$ws.Cells.Item(68, 11).Value = 2070  # K68: 2205 -> 2070
$ws.Cells.Item(68, 12).Value = 63303  # L68: 53153.50199999999 -> 63303
$ws.Cells.Item(68, 13).Value = -1259  # M68: -1394 -> -1259
$ws.Cells.Item(68, 14).Value = -64925  # N68: -54775.50199999999 -> -64925

$ws.Cells.Item(71, 8).Value = 13446.875  # H71: 13472.125 -> 13446.875
$ws.Cells.Item(71, 9).Value = 690  # I71: 735 -> 690
$ws.Cells.Item(71, 10).Value = 21101  # J71: 17717.834 -> 21101
$ws.Cells.Item(71, 11).Value = 6210  # K71: 6615 -> 6210
$ws.Cells.Item(71, 12).Value = 189909  # L71: 159460.506 -> 189909
$ws.Cells.Item(71, 13).Value = -2154  # M71: -2559 -> -2154
$ws.Cells.Item(71, 14).Value = -198021  # N71: -167572.506 -> -198021

$ws.Cells.Item(107, 8).Value = 24842  # H107: 24254.857 -> 24842
$ws.Cells.Item(107, 9).Value = 465.85715  # I107: 480.07693 -> 465.85715
$ws.Cells.Item(107, 10).Value = 37481.48  # J107: 34912.516 -> 37481.48
$ws.Cells.Item(107, 11).Value = 1397.57145  # K107: 1440.23079 -> 1397.57145
$ws.Cells.Item(107, 12).Value = 112444.44  # L107: 104737.548 -> 112444.44
$ws.Cells.Item(107, 13).Value = 522.4285500000001  # M107: 479.7692099999999 -> 522.4285500000001
$ws.Cells.Item(107, 14).Value = -116284.44  # N107: -108577.548 -> -116284.44

$ws.Cells.Item(122, 8).Value = 3066.257  # H122: 2381.4565 -> 3066.257
$ws.Cells.Item(122, 9).Value = 778.5  # I122: 471.73334 -> 778.5
$ws.Cells.Item(122, 10).Value = 3744.111  # J122: 3305.516 -> 3744.111
$ws.Cells.Item(122, 11).Value = 7006.5  # K122: 4245.60006 -> 7006.5
$ws.Cells.Item(122, 12).Value = 33696.999  # L122: 29749.644 -> 33696.999
$ws.Cells.Item(122, 13).Value = -4556.5  # M122: -1795.60006 -> -4556.5
$ws.Cells.Item(122, 14).Value = -38596.999  # N122: -34649.644 -> -38596.999

$ws.Cells.Item(131, 8).Value = 875.3605  # H131: 856.8936 -> 875.3605
$ws.Cells.Item(131, 9).Value = 584.5  # I131: 427.22223 -> 584.5
$ws.Cells.Item(131, 10).Value = 913.6316  # J131: 902.38824 -> 913.6316
$ws.Cells.Item(131, 11).Value = 1753.5  # K131: 1281.66669 -> 1753.5
$ws.Cells.Item(131, 12).Value = 2740.8948  # L131: 2707.16472 -> 2740.8948
$ws.Cells.Item(131, 13).Value = 3286.5  # M131: 3758.33331 -> 3286.5
$ws.Cells.Item(131, 14).Value = -12820.8948  # N131: -12787.16472 -> -12820.8948

$ws.Cells.Item(132, 8).Value = 1610.6316  # H132: 1880.1333 -> 1610.6316
$ws.Cells.Item(132, 9).Value = 698.4545000000001  # I132: 754.7143 -> 698.4545000000001
$ws.Cells.Item(132, 11).Value = 6286.0905  # K132: 6792.428699999999 -> 6286.0905
$ws.Cells.Item(132, 13).Value = -3756.0905  # M132: -4262.428699999999 -> -3756.0905

$ws.Cells.Item(133, 8).Value = 3379.4736  # H133: 3369.5 -> 3379.4736
$ws.Cells.Item(133, 9).Value = 4632.857  # I133: 4732.857 -> 4632.857
$ws.Cells.Item(133, 10).Value = 2648.3333  # J133: 2635.3845 -> 2648.3333
$ws.Cells.Item(133, 11).Value = 13898.571  # K133: 14198.571 -> 13898.571
$ws.Cells.Item(133, 12).Value = 7944.999899999999  # L133: 7906.1535 -> 7944.999899999999
$ws.Cells.Item(133, 13).Value = -8838.571  # M133: -9138.571 -> -8838.571
$ws.Cells.Item(133, 14).Value = -18064.9999  # N133: -18026.1535 -> -18064.9999

$ws.Cells.Item(135, 8).Value = 1375.0526  # H135: 1435.4722 -> 1375.0526
$ws.Cells.Item(135, 9).Value = 328.6316  # I135: 342.72223 -> 328.6316
$ws.Cells.Item(135, 10).Value = 2421.4736  # J135: 2528.2222 -> 2421.4736
$ws.Cells.Item(135, 11).Value = 2957.6844  # K135: 3084.50007 -> 2957.6844
$ws.Cells.Item(135, 12).Value = 21793.2624  # L135: 22753.9998 -> 21793.2624
$ws.Cells.Item(135, 13).Value = -422.6844000000001  # M135: -549.5000700000001 -> -422.6844000000001
$ws.Cells.Item(135, 14).Value = -26863.2624  # N135: -27823.9998 -> -26863.2624

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 736.0714  # H97: 891.1111 -> 736.0714
$ws.Cells.Item(97, 9).Value = 761.9231  # I97: 915 -> 761.9231
$ws.Cells.Item(97, 10).Value = 400  # J97: 700 -> 400
$ws.Cells.Item(97, 11).Value = 761.9231  # K97: 915 -> 761.9231
$ws.Cells.Item(97, 12).Value = 400  # L97: 700 -> 400
$ws.Cells.Item(97, 13).Value = -265.9231  # M97: -419 -> -265.9231
$ws.Cells.Item(97, 14).Value = -1392  # N97: -1692 -> -1392

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 8534.549999999999  # H40: 7694.8 -> 8534.549999999999
$ws.Cells.Item(40, 9).Value = 6401.143  # I40: 5223.6665 -> 6401.143
$ws.Cells.Item(40, 10).Value = 9683.308000000001  # J40: 9716.637000000001 -> 9683.308000000001
$ws.Cells.Item(40, 11).Value = 6401.143  # K40: 5223.6665 -> 6401.143
$ws.Cells.Item(40, 12).Value = 9683.308000000001  # L40: 9716.637000000001 -> 9683.308000000001
$ws.Cells.Item(40, 13).Value = -6265.143  # M40: -5087.6665 -> -6265.143
$ws.Cells.Item(40, 14).Value = -9955.308000000001  # N40: -9988.637000000001 -> -9955.308000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5556853  # H132: 6062027 -> 5556853
$ws.Cells.Item(132, 9).Value = 596.7027  # I132: 657.9697 -> 596.7027
$ws.Cells.Item(132, 10).Value = 14495178  # J132: 15154080 -> 14495178
$ws.Cells.Item(132, 11).Value = 1790.1081  # K132: 1973.9091 -> 1790.1081
$ws.Cells.Item(132, 12).Value = 43485534  # L132: 45462240 -> 43485534
$ws.Cells.Item(132, 13).Value = 739.8918999999999  # M132: 556.0909000000001 -> 739.8918999999999
$ws.Cells.Item(132, 14).Value = -43490594  # N132: -45467300 -> -43490594

$ws.Cells.Item(136, 8).Value = 1947.7656  # H136: 2084.1018 -> 1947.7656
$ws.Cells.Item(136, 9).Value = 712.9761999999999  # I136: 763.5135 -> 712.9761999999999
$ws.Cells.Item(136, 11).Value = 2138.9286  # K136: 2290.5405 -> 2138.9286
$ws.Cells.Item(136, 13).Value = 411.0714000000003  # M136: 259.4594999999999 -> 411.0714000000003
